# "ran final round of usda families resazurin data"
#
# Appends the 2025-06-10 timepoint/temperature readings (rows 36:42) to the
# resazurin assay log on Sheet1, following the same layout used for every
# previous sampling date: column A = date (yyyymmdd), column B = timepoint
# (0-6 hours), column C = temperature-group label ("18C" for timepoints
# 0-2, "42C" for timepoints 3-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$date = 20250610
$labels = @("18C", "18C", "18C", "42C", "42C", "42C", "42C")

for ($t = 0; $t -le 6; $t++) {
    $row = 36 + $t
    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $t
    $ws.Cells.Item($row, 3).Value = $labels[$t]
}

# Column C had no formatting of its own yet for these rows; pick up the
# same cell style already used by column B (and by every earlier date
# block) so the new temperature labels render consistently.
$ws.Range("B36:B42").Copy() | Out-Null
$ws.Range("C36:C42").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Match the author's final on-screen state: scrolled down with B40 selected.
$ws.Range("B40").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
